$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D28").Value = "Data Driven Control :: 모델 학습 (Linear, Bayesian, Gaussian Process Regression)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/252"

$ws.Range("D51").Value = "[css] @font-face로 폰트를 등록할 때 굵기별로 지정하는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/unicode-range%EB%A1%9C-%ED%8A%B9%EC%A0%95-%ED%8F%B0%ED%8A%B8-%ED%8A%B9%EC%A0%95-%EC%96%B8%EC%96%B4%EC%97%90%EB%A7%8C-%EC%A0%81%EC%9A%A9%EB%90%98%EA%B2%8C-%ED%95%98%EA%B8%B0"
